$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.089.72"
$ws.Range("E2").Value = "  +1.04%  "

$ws.Range("D3").Value = "2.357.93"
$ws.Range("E3").Value = "  +0.74%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.17%  "

$ws.Range("D5").Value = "'551.18"
$ws.Range("E5").Value = "  +1.32%  "

$ws.Range("D6").Value = "'132.59"
$ws.Range("E6").Value = "  -1.31%  "

$ws.Range("E7").Value = "  -0.17%  "

$ws.Range("D8").Value = "'0.569"
$ws.Range("E8").Value = "  +0.09%  "

$ws.Range("E9").Value = "  +4.20%  "

$ws.Range("E10").Value = "  +4.69%  "

$ws.Range("E11").Value = "  -1.16%  "

$ws.Range("D12").Value = "'0.354"
$ws.Range("E12").Value = "  -1.09%  "

$ws.Range("D13").Value = "'24.08"
$ws.Range("E13").Value = "  +2.45%  "

$ws.Range("D14").Value = "2.776.50"
$ws.Range("E14").Value = "  +0.66%  "

$ws.Range("D15").Value = "57.967.00"
$ws.Range("E15").Value = "  +0.82%  "

$ws.Range("E16").Value = "  +2.41%  "

$ws.Range("D17").Value = "2.341.05"
$ws.Range("E17").Value = "  -0.28%  "

$ws.Range("D18").Value = "'10.99"
$ws.Range("E18").Value = "  +3.70%  "

$ws.Range("D19").Value = "'4.31"
$ws.Range("E19").Value = "  +2.02%  "

$ws.Range("D20").Value = "'330.47"
$ws.Range("E20").Value = "  -0.91%  "

$ws.Range("D21").Value = "'6.88"
$ws.Range("E21").Value = "  +2.83%  "

$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = "  -0.02%  "

$ws.Range("D23").Value = "'64.03"
$ws.Range("E23").Value = "  +3.19%  "

$ws.Range("E24").Value = "  +0.66%  "

$ws.Range("D25").Value = "'0.998"
$ws.Range("E25").Value = "  -0.49%  "

$ws.Range("D26").Value = "'8.27"
$ws.Range("E26").Value = "  -2.40%  "

$ws.Range("E27").Value = "  -4.69%  "

$ws.Range("E28").Value = "  +0.10%  "

$ws.Range("D29").Value = "'170.90"
$ws.Range("E29").Value = "  +0.47%  "

$ws.Range("D30").Value = "0.0₃0738"
$ws.Range("E30").Value = "  +0.89%  "

$ws.Range("D31").Value = "'6.14"
$ws.Range("E31").Value = "  +0.17%  "

$ws.Range("D32").Value = "'18.40"
$ws.Range("E32").Value = "  -0.41%  "

$ws.Range("E33").Value = "  -2.48%  "

$ws.Range("E34").Value = "  -0.06%  "

$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.26%  "

$ws.Range("D36").Value = "'4.15"
$ws.Range("E36").Value = "  -0.08%  "

$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "'1.24"
$ws.Range("E37").Value = "  -0.73%  "

$ws.Range("B38").Value = "PolygonEcosystemToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D38").Value = "'0.432"
$ws.Range("E38").Value = "  +14.72%  "

$ws.Range("D39").Value = "'40.45"
$ws.Range("E39").Value = "  +3.74%  "

$ws.Range("E40").Value = "  -0.69%  "

$ws.Range("D41").Value = "'141.93"
$ws.Range("E41").Value = "  -3.84%  "

$ws.Range("D42").Value = "'3.67"
$ws.Range("E42").Value = "  +1.39%  "

$ws.Range("D43").Value = "'289.57"
$ws.Range("E43").Value = "  +2.41%  "

$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D44").Value = "'0.0955"
$ws.Range("E44").Value = "  +1.63%  "

$ws.Range("B45").Value = "Polygon"
$ws.Range("C45").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D45").Value = "'0.412"
$ws.Range("E45").Value = "  +7.02%  "

$ws.Range("E46").Value = "  +1.95%  "

$ws.Range("D47").Value = "'0.568"
$ws.Range("E47").Value = "  +1.71%  "

$ws.Range("E48").Value = "  -1.96%  "

$ws.Range("E49").Value = "  +2.24%  "

$ws.Range("E50").Value = "  -0.04%  "

$ws.Range("E51").Value = "  -0.04%  "

